# Apply the Jan 5 2023 GitHub Actions crypto-symbol-list refresh.
# Price (D) and Volume(1h) (E) columns are stored as literal text in the
# workbook (e.g. "256.91", "0.44%"), so numeric-looking values are written
# with a leading apostrophe to force Excel to keep them as text instead of
# auto-converting them to numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''256.91'
$ws.Range("E2").Value = '''0.44%'

$ws.Range("D3").Value = '''27.04'
$ws.Range("E3").Value = '''-3.82%'

$ws.Range("D4").Value = '''4.629'
$ws.Range("E4").Value = '''-11.31%'

$ws.Range("D5").Value = '''0.05894'
$ws.Range("E5").Value = '''0.45%'

$ws.Range("D6").Value = '''6.636'
$ws.Range("E6").Value = '''-0.85%'

$ws.Range("D7").Value = '''0.8650'
$ws.Range("E7").Value = '''-0.55%'

$ws.Range("E8").Value = '''-2.65%'

$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = '''0.01045'
$ws.Range("E9").Value = '''1,626.95%'

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1405'
$ws.Range("E10").Value = '''-0.38%'

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.03844'
$ws.Range("E11").Value = '''10.10%'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.07077'
$ws.Range("E12").Value = '''-0.54%'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03201'
$ws.Range("E13").Value = '''0.96%'

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09249'
$ws.Range("E14").Value = '''0.24%'

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001555'
$ws.Range("E15").Value = '''0.03%'

$ws.Range("D16").Value = '''0.006085'
$ws.Range("E16").Value = '''0.93%'

$ws.Range("D17").Value = '''3.515'
$ws.Range("E17").Value = '''0.50%'

$ws.Range("D18").Value = '''3.187'
$ws.Range("E18").Value = '''-0.81%'

$ws.Range("E19").Value = '''-0.60%'

$ws.Range("D20").Value = '''0.3098'
$ws.Range("E20").Value = '''-2.39%'

$ws.Range("E21").Value = '''-1.76%'

$ws.Range("E22").Value = '''9.16%'

$ws.Range("D23").Value = '''0.04224'
$ws.Range("E23").Value = '''1.28%'

$ws.Range("D24").Value = '''0.001217'
$ws.Range("E24").Value = '''-0.75%'

$ws.Range("D25").Value = '''0.004281'
$ws.Range("E25").Value = '''-6.06%'

$ws.Range("D26").Value = '''0.0001198'
$ws.Range("E26").Value = '''-0.12%'

$ws.Range("D27").Value = '''0.0001934'
$ws.Range("E27").Value = '''31.94%'

$ws.Range("D40").Value = '''0.03829'
$ws.Range("E40").Value = '''0.18%'

$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '''0.006132'
$ws.Range("E41").Value = '''11.10%'

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1100'
$ws.Range("E42").Value = '''-0.19%'

$ws.Range("D43").Value = '''0.002311'
$ws.Range("E43").Value = '''-1.41%'

$ws.Range("D44").Value = '''0.01162'
$ws.Range("E44").Value = '''21.58%'

$ws.Range("D45").Value = '''0.00005462'
$ws.Range("E45").Value = '''1.30%'

$ws.Range("E46").Value = '''-0.12%'

$ws.Range("D47").Value = '''0.07770'
$ws.Range("E47").Value = '''-18.20%'

$ws.Range("E48").Value = '''6.85%'

$ws.Range("D49").Value = '''0.00002097'
$ws.Range("E49").Value = '''-0.12%'

$ws.Range("D50").Value = '''0.0001997'
$ws.Range("E50").Value = '''-0.12%'
